# Add two new columns (I: "I0", J: "IF") to the sheet, mirroring the
# existing H column's header style, and populate data rows 2-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell H1 (which already carries the bold/centered/
# bordered header style) into I1 and J1 so the new headers pick up the
# same cell style (s="1") as the other header cells, then overwrite the
# copied text with the correct header labels.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for the new I and J columns, one pair per data row (rows 2-24).
$ijValues = @(
    @(7, 8),
    @(4, 4),
    @(6, 7),
    @(6, 7),
    @(7, 8),
    @(8, 8),
    @(10, 11),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(4, 4),
    @(7, 7),
    @(6, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(4, 4),
    @(7, 7),
    @(9, 9)
)

$row = 2
foreach ($pair in $ijValues) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
